$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add three new "Ice corer" gear-type rows to the alphabetically-sorted
# list, inserting each just above "Ice corer 14 cm" (originally row 8)
# so the new rows inherit formatting from their neighbours, matching how
# the rows were entered by hand.

# Insert "Ice corer 18 cm" directly above "Ice corer 14 cm" (row 8).
$ws.Rows("8:8").Insert()
$ws.Range("A8").Value = "Ice corer 18 cm"

# Insert "Ice corer 12 cm" directly below "Ice corer 14 cm" (now row 9).
$ws.Rows("10:10").Insert()
$ws.Range("A10").Value = "Ice corer 12 cm"

# Insert "Ice corer 22 cm" above "Ice corer 18 cm" (row 8).
$ws.Rows("8:8").Insert()
$ws.Range("A8").Value = "Ice corer 22 cm"

$ws.Range("A8").Select()
